$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91828
$ws.Range("B3").Value = 57881
$ws.Range("B4").Value = 79243
$ws.Range("B5").Value = 91828
$ws.Range("B6").Value = 57988
$ws.Range("B7").Value = 79243
$ws.Range("B8").Value = 80349
